$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 132.83333
$ws.Cells.Item(4, 9).Value = 119.4
$ws.Cells.Item(4, 11).Value = 119.4
$ws.Cells.Item(4, 13).Value = -5.400000000000006
$ws.Cells.Item(33, 8).Value = 346.05884
$ws.Cells.Item(33, 9).Value = 355.7857
$ws.Cells.Item(33, 11).Value = 355.7857
$ws.Cells.Item(33, 13).Value = -126.7857
$ws.Cells.Item(40, 8).Value = 38166.168
$ws.Cells.Item(40, 9).Value = 52500
$ws.Cells.Item(40, 10).Value = 30999.25
$ws.Cells.Item(40, 11).Value = 52500
$ws.Cells.Item(40, 12).Value = 30999.25
$ws.Cells.Item(40, 13).Value = -52325
$ws.Cells.Item(40, 14).Value = -31349.25
$ws.Cells.Item(41, 8).Value = 1533.625
$ws.Cells.Item(41, 10).Value = 2150.6365
$ws.Cells.Item(41, 12).Value = 2150.6365
$ws.Cells.Item(41, 14).Value = -3030.6365
$ws.Cells.Item(74, 8).Value = 5099.25
$ws.Cells.Item(74, 9).Value = 5099.25
$ws.Cells.Item(74, 11).Value = 5099.25
$ws.Cells.Item(74, 13).Value = -4163.25
$ws.Cells.Item(77, 8).Value = 5099.25
$ws.Cells.Item(77, 9).Value = 5099.25
$ws.Cells.Item(77, 11).Value = 25496.25
$ws.Cells.Item(77, 13).Value = -20816.25
$ws.Cells.Item(132, 8).Value = 8536.773999999999
$ws.Cells.Item(132, 9).Value = 1871.5
$ws.Cells.Item(132, 11).Value = 5614.5
$ws.Cells.Item(132, 13).Value = -3084.5
$ws.Cells.Item(136, 8).Value = 109853
$ws.Cells.Item(136, 10).Value = 109853
$ws.Cells.Item(136, 12).Value = 109853
$ws.Cells.Item(136, 14).Value = -120053
$ws.Cells.Item(137, 8).Value = 12348942
$ws.Cells.Item(137, 9).Value = 2072.8333
$ws.Cells.Item(137, 11).Value = 6218.499899999999
$ws.Cells.Item(137, 13).Value = -3668.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 958.3043
$ws.Cells.Item(74, 9).Value = 768.94116
$ws.Cells.Item(74, 11).Value = 768.94116
$ws.Cells.Item(74, 13).Value = 105.05884
$ws.Cells.Item(77, 8).Value = 958.3043
$ws.Cells.Item(77, 9).Value = 768.94116
$ws.Cells.Item(77, 11).Value = 3844.7058
$ws.Cells.Item(77, 13).Value = 523.2942000000003
$ws.Cells.Item(103, 8).Value = 42500
$ws.Cells.Item(103, 10).Value = 42500
$ws.Cells.Item(103, 12).Value = 42500
$ws.Cells.Item(103, 14).Value = -44844
$ws.Cells.Item(122, 8).Value = 7501.5
$ws.Cells.Item(122, 9).Value = 6801.5454
$ws.Cells.Item(122, 11).Value = 20404.6362
$ws.Cells.Item(122, 13).Value = -17954.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3144.95
$ws.Cells.Item(86, 9).Value = 1872.1428
$ws.Cells.Item(86, 11).Value = 1872.1428
$ws.Cells.Item(86, 13).Value = -749.1428000000001
$ws.Cells.Item(89, 8).Value = 3144.95
$ws.Cells.Item(89, 9).Value = 1872.1428
$ws.Cells.Item(89, 11).Value = 9360.714
$ws.Cells.Item(89, 13).Value = -3744.714
$ws.Cells.Item(107, 8).Value = 1424.1666
$ws.Cells.Item(107, 9).Value = 901.5714
$ws.Cells.Item(107, 11).Value = 901.5714
$ws.Cells.Item(107, 13).Value = 1018.4286
$ws.Cells.Item(111, 8).Value = 75000
$ws.Cells.Item(111, 10).Value = 75000
$ws.Cells.Item(111, 12).Value = 75000
$ws.Cells.Item(111, 14).Value = -83180

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1960
$ws.Cells.Item(16, 9).Value = 1651.8572
$ws.Cells.Item(16, 11).Value = 1651.8572
$ws.Cells.Item(16, 13).Value = -1364.8572
$ws.Cells.Item(22, 8).Value = 1894.4286
$ws.Cells.Item(22, 9).Value = 1757.3334
$ws.Cells.Item(22, 10).Value = 1997.25
$ws.Cells.Item(22, 11).Value = 1757.3334
$ws.Cells.Item(22, 12).Value = 1997.25
$ws.Cells.Item(22, 13).Value = -1407.3334
$ws.Cells.Item(22, 14).Value = -2697.25
$ws.Cells.Item(31, 8).Value = 11365645
$ws.Cells.Item(31, 9).Value = 12988576
$ws.Cells.Item(31, 11).Value = 12988576
$ws.Cells.Item(31, 13).Value = -12988281
$ws.Cells.Item(34, 8).Value = 11365645
$ws.Cells.Item(34, 9).Value = 12988576
$ws.Cells.Item(34, 11).Value = 12988576
$ws.Cells.Item(34, 13).Value = -12988374
$ws.Cells.Item(107, 8).Value = 909835.25
$ws.Cells.Item(107, 9).Value = 2597934.5
$ws.Cells.Item(107, 10).Value = 858.7692
$ws.Cells.Item(107, 11).Value = 2597934.5
$ws.Cells.Item(107, 12).Value = 858.7692
$ws.Cells.Item(107, 13).Value = -2596014.5
$ws.Cells.Item(107, 14).Value = -4698.7692
$ws.Cells.Item(109, 8).Value = 43250
$ws.Cells.Item(109, 10).Value = 43250
$ws.Cells.Item(109, 12).Value = 43250
$ws.Cells.Item(109, 14).Value = -45330
$ws.Cells.Item(113, 8).Value = 1960
$ws.Cells.Item(113, 9).Value = 1651.8572
$ws.Cells.Item(113, 11).Value = 1651.8572
$ws.Cells.Item(113, 13).Value = 518.1428000000001
$ws.Cells.Item(132, 8).Value = 70183350
$ws.Cells.Item(132, 9).Value = 95240200
$ws.Cells.Item(132, 11).Value = 285720600
$ws.Cells.Item(132, 13).Value = -285718070
$ws.Cells.Item(134, 8).Value = 3311.8125
$ws.Cells.Item(134, 9).Value = 2621.5557
$ws.Cells.Item(134, 11).Value = 7864.6671
$ws.Cells.Item(134, 13).Value = -5329.6671
$ws.Cells.Item(141, 8).Value = 118057.02
$ws.Cells.Item(141, 10).Value = 118986.66
$ws.Cells.Item(141, 12).Value = 118986.66
$ws.Cells.Item(141, 14).Value = -129346.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(28, 8).Value = 2050
$ws.Cells.Item(28, 9).Value = 1700
$ws.Cells.Item(28, 11).Value = 5100
$ws.Cells.Item(28, 13).Value = -4868
$ws.Cells.Item(37, 8).Value = 333413340
$ws.Cells.Item(37, 10).Value = 333413340
$ws.Cells.Item(37, 12).Value = 1000240020
$ws.Cells.Item(37, 14).Value = -1000240244
$ws.Cells.Item(68, 8).Value = 1112.091
$ws.Cells.Item(68, 10).Value = 1112.091
$ws.Cells.Item(68, 12).Value = 3336.273
$ws.Cells.Item(68, 14).Value = -4958.272999999999
$ws.Cells.Item(71, 8).Value = 1112.091
$ws.Cells.Item(71, 10).Value = 1112.091
$ws.Cells.Item(71, 12).Value = 10008.819
$ws.Cells.Item(71, 14).Value = -18120.819
$ws.Cells.Item(80, 8).Value = 5416.5
$ws.Cells.Item(80, 10).Value = 5416.5
$ws.Cells.Item(80, 12).Value = 16249.5
$ws.Cells.Item(80, 14).Value = -18121.5
$ws.Cells.Item(83, 8).Value = 5416.5
$ws.Cells.Item(83, 10).Value = 5416.5
$ws.Cells.Item(83, 12).Value = 48748.5
$ws.Cells.Item(83, 14).Value = -58108.5
$ws.Cells.Item(107, 8).Value = 395.8889
$ws.Cells.Item(107, 10).Value = 483.16666
$ws.Cells.Item(107, 12).Value = 1449.49998
$ws.Cells.Item(107, 14).Value = -5289.499980000001
$ws.Cells.Item(115, 8).Value = 6162.294
$ws.Cells.Item(115, 9).Value = 4581
$ws.Cells.Item(115, 10).Value = 7024.8184
$ws.Cells.Item(115, 11).Value = 13743
$ws.Cells.Item(115, 12).Value = 21074.4552
$ws.Cells.Item(115, 13).Value = -12568
$ws.Cells.Item(115, 14).Value = -23424.4552
$ws.Cells.Item(128, 8).Value = 118474.25
$ws.Cells.Item(128, 9).Value = 118474.25
$ws.Cells.Item(128, 11).Value = 355422.75
$ws.Cells.Item(128, 13).Value = -350442.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 980006.75
$ws.Cells.Item(80, 9).Value = 1387459.6
$ws.Cells.Item(80, 11).Value = 1387459.6
$ws.Cells.Item(80, 13).Value = -1386461.6
$ws.Cells.Item(83, 8).Value = 980006.75
$ws.Cells.Item(83, 9).Value = 1387459.6
$ws.Cells.Item(83, 11).Value = 6937298
$ws.Cells.Item(83, 13).Value = -6932306
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 404560.62
$ws.Cells.Item(122, 9).Value = 771540.0600000001
$ws.Cells.Item(122, 11).Value = 2314620.18
$ws.Cells.Item(122, 13).Value = -2312170.18
$ws.Cells.Item(123, 8).Value = 39243.75
$ws.Cells.Item(123, 10).Value = 39243.75
$ws.Cells.Item(123, 12).Value = 39243.75
$ws.Cells.Item(123, 14).Value = -44143.75
$ws.Cells.Item(134, 8).Value = 543227.8
$ws.Cells.Item(134, 10).Value = 543227.8
$ws.Cells.Item(134, 12).Value = 1629683.4
$ws.Cells.Item(134, 14).Value = -1634753.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1230.1875
$ws.Cells.Item(16, 9).Value = 656
$ws.Cells.Item(16, 10).Value = 5249.5
$ws.Cells.Item(16, 11).Value = 656
$ws.Cells.Item(16, 12).Value = 5249.5
$ws.Cells.Item(16, 13).Value = -486
$ws.Cells.Item(16, 14).Value = -5589.5
$ws.Cells.Item(40, 8).Value = 7816235.5
$ws.Cells.Item(40, 9).Value = 10419899
$ws.Cells.Item(40, 10).Value = 5246
$ws.Cells.Item(40, 11).Value = 10419899
$ws.Cells.Item(40, 12).Value = 5246
$ws.Cells.Item(40, 13).Value = -10419763
$ws.Cells.Item(40, 14).Value = -5518
$ws.Cells.Item(82, 8).Value = 1564215.4
$ws.Cells.Item(82, 9).Value = 2605227.8
$ws.Cells.Item(82, 11).Value = 2605227.8
$ws.Cells.Item(82, 13).Value = -2604866.8
$ws.Cells.Item(85, 8).Value = 1564215.4
$ws.Cells.Item(85, 9).Value = 2605227.8
$ws.Cells.Item(85, 11).Value = 2605227.8
$ws.Cells.Item(85, 13).Value = -2603979.8
$ws.Cells.Item(100, 8).Value = 1799.8
$ws.Cells.Item(100, 9).Value = 1666.6666
$ws.Cells.Item(100, 10).Value = 1999.5
$ws.Cells.Item(100, 11).Value = 1666.6666
$ws.Cells.Item(100, 12).Value = 1999.5
$ws.Cells.Item(100, 13).Value = -1125.6666
$ws.Cells.Item(100, 14).Value = -3081.5
$ws.Cells.Item(122, 8).Value = 5646.0625
$ws.Cells.Item(122, 9).Value = 3082.5715
$ws.Cells.Item(122, 10).Value = 8739.931
$ws.Cells.Item(122, 11).Value = 9247.7145
$ws.Cells.Item(122, 12).Value = 26219.793
$ws.Cells.Item(122, 13).Value = -6797.7145
$ws.Cells.Item(122, 14).Value = -31119.793
$ws.Cells.Item(136, 8).Value = 3057.3208
$ws.Cells.Item(136, 9).Value = 2716.617
$ws.Cells.Item(136, 11).Value = 8149.851000000001
$ws.Cells.Item(136, 13).Value = -5599.851000000001
$ws.Cells.Item(139, 8).Value = 83657
$ws.Cells.Item(139, 10).Value = 83657
$ws.Cells.Item(139, 12).Value = 83657
$ws.Cells.Item(139, 14).Value = -93937

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2160.647
$ws.Cells.Item(136, 9).Value = 2083.2188
$ws.Cells.Item(136, 11).Value = 6249.6564
$ws.Cells.Item(136, 13).Value = -3699.6564
